$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("ECs","Jam2","Jam3","ECs",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,1429.112351066053,5716.449404264211,0.1094779631300779,0.07900592227825659)
    ,@("ECs","Jam2","Jam3","FAPs",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,953.0761855138794,5718.457113083276,0.07301094236573515,0.07903367042674858)
    ,@("ECs","Jam2","Jam3","M1",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,64.610139187702,387.66083512621196,0.004949496399316693,0.005357784114639207)
    ,@("ECs","Jam2","Jam3","M2",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,109.14284612047952,654.857076722877,0.008360949700404173,0.009050650788291774)
    ,@("ECs","Jam2","Jam3","Neutro",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,1097.8872668565587,6587.32360113935,0.08410427747842548,0.09104210317423925)
    ,@("ECs","Jam2","Jam3","sCs",2,1,93.83570850000001,187.671417,0.5679500693534724,0.4713618706519805,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,3760.1240808157254,15040.496323262902,0.288046440279513,0.20787173986980512)
    ,@("FAPs","Jam2","Jam3","ECs",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,53.889533623951,323.33720174370603,0.004128238322745641,0.004468779835884911)
    ,@("FAPs","Jam2","Jam3","FAPs",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,35.938973662304655,323.450762960742,0.0027531254842207862,0.004470349343117892)
    ,@("FAPs","Jam2","Jam3","M1",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,2.4363446761946665,21.927102085752,0.0001866375673218665,0.0003030501629004275)
    ,@("FAPs","Jam2","Jam3","M2",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,4.115601598038,37.040414382342,0.0003152779973330934,0.0005119282780081566)
    ,@("FAPs","Jam2","Jam3","Neutro",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,41.399567177794,372.596104600146,0.0031714373511017493,0.0051495774386215)
    ,@("FAPs","Jam2","Jam3","sCs",3,1,3.538394,10.615182,0.02141648579016069,0.02666144970190763,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,141.78824564226395,850.729473853584,0.010861769067437552,0.011757764643374745)
    ,@("M1","Jam2","Jam3","ECs",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,5.158259346518,30.949556079108,0.00039515138619593146,0.00042774772401703466)
    ,@("M1","Jam2","Jam3","FAPs",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,3.4400473400173324,30.960426060156,0.0002635267798051055,0.00042789795588599926)
    ,@("M1","Jam2","Jam3","M1",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,0.23320479603733332,2.098843164336,0.000017864785818475165,0.000029007698343675563)
    ,@("M1","Jam2","Jam3","M2",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,0.393941809884,3.545476288956,0.000030178136033675184,0.0000490013300766219)
    ,@("M1","Jam2","Jam3","Neutro",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,3.962730607892,35.664575471027995,0.0003035672283299581,0.0004929130791657436)
    ,@("M1","Jam2","Jam3","sCs",3,1,0.338692,1.016076,0.002049967416076645,0.002552010805591038,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,13.571847706351996,81.43108623811199,0.0010396790998934997,0.0011254430181019631)
    ,@("M2","Jam2","Jam3","ECs",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,3.4768788815651663,20.861273289391,0.00026634828095901383,0.00028831955284999814)
    ,@("M2","Jam2","Jam3","FAPs",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,2.3187333448374434,20.868600103537,0.00017762788450133776,0.0002884208153064688)
    ,@("M2","Jam2","Jam3","M1",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,0.15718962075244441,1.414706586772,0.000012041600152941135,0.000019552381336161363)
    ,@("M2","Jam2","Jam3","M2",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,0.265532976793,2.389796791137,0.000020341304461807346,0.00003302891116302972)
    ,@("M2","Jam2","Jam3","Neutro",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,2.671043357525666,24.039390217730997,0.00020461679307151896,0.00033224368149606636)
    ,@("M2","Jam2","Jam3","sCs",3,1,0.2282923333333333,0.684877,0.001381762322917109,0.001720160208981192,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,9.147983361070661,54.88790016642399,0.0007007864597704899,0.0007585948668294678)
    ,@("Neutro","Jam2","Jam3","ECs",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,968.7060184960303,5812.236110976183,0.07420827459624363,0.08032977150189023)
    ,@("Neutro","Jam2","Jam3","FAPs",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,646.0308290694087,5814.277461624681,0.04948955848922301,0.08035798460748644)
    ,@("Neutro","Jam2","Jam3","M1",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,43.795178622803995,394.156607605236,0.0033549545261196247,0.005447560907771003)
    ,@("Neutro","Jam2","Jam3","M2",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,73.98111970260899,665.830077323481,0.005667365682678758,0.009202306470219771)
    ,@("Neutro","Jam2","Jam3","Neutro",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,744.1892180420668,6697.702962378603,0.05700903761263855,0.09256763460441104)
    ,@("Neutro","Jam2","Jam3","sCs",3,1,63.60536699999999,190.816101,0.3849778850329995,0.4792601652167268,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,2548.7533045676505,15292.51982740591,0.19524869412609594,0.21135490712494837)
    ,@("sCs","Jam2","Jam3","ECs",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,2,1,15.2299415,30.459883,0.1927598375940027,0.1676120348236416,55.921025060753,223.684100243012,0.004283861877780607,0.0030914939307428112)
    ,@("sCs","Jam2","Jam3","FAPs",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,3,1,10.15686033333333,30.470581,0.1285516919627237,0.1676709028615964,37.293776948447324,223.76266169068398,0.002856910959238334,0.0030925797130510056)
    ,@("sCs","Jam2","Jam3","M1",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,3,1,0.6885453333333333,2.065636,0.008714668183685525,0.01136660482789667,2.5281883611173335,15.169130166703999,0.00019367330495592564,0.00020964956290619557)
    ,@("sCs","Jam2","Jam3","M2",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,2,0.6666666666666666,1.163127,3.489381,0.01472127595639153,0.01920106684864657,4.270748782314,25.624492693883997,0.0003271631354800229,0.00035415107088721513)
    ,@("sCs","Jam2","Jam3","Neutro",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,3,1,11.700101,35.100303,0.1480839285294319,0.1931469404776233,42.960220249982,257.76132149989195,0.0032909920658646496,0.003562468499689695)
    ,@("sCs","Jam2","Jam3","sCs",2,1,3.671782,7.343564,0.02222383008437381,0.01844434341481282,2,1,40.07135599999999,80.14271199999999,0.5071685977737647,0.4410024501605956,147.13328367639195,588.5331347055679,0.011271228741054271,0.008134000637535901)
)

$numRows = $rows.Count
$numCols = 20

$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T37").Value = $data
